$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.405.88"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "3.580.67"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'207.19"
$ws.Range("E5").Value = "  +8.56%  "
$ws.Range("D6").Value = "'567.98"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.610"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.677"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'63.63"
$ws.Range("E10").Value = "  +14.70%  "
$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "'0.0000282"
$ws.Range("E12").Value = "  +5.04%  "
$ws.Range("D13").Value = "'10.13"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").Value = "4.160.33"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "3.611.59"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "'19.17"
$ws.Range("E16").Value = "  +5.45%  "
$ws.Range("D17").Value = "'0.125"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "68.185.87"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "'12.19"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'1.06"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "'404.97"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("D24").Value = "'84.70"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").Value = "'2.88"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").Value = "'12.45"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'3.81"
$ws.Range("E27").Value = "  +4.60%  "
$ws.Range("D28").Value = "'9.08"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "'734.90"
$ws.Range("E29").Value = "  +15.24%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'31.47"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.55"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "'12.06"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "'63.94"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'0.112"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "'41.30"
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("E36").Value = "  +4.93%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  +8.29%  "
$ws.Range("D39").Value = "'3.12"
$ws.Range("E39").Value = "  +27.15%  "
$ws.Range("D40").Value = "0.0₃0736"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "3.165.82"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'2.58"
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("D45").Value = "'2.75"
$ws.Range("E45").Value = "  +9.36%  "
$ws.Range("D46").Value = "'0.0412"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'8.77"
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("D48").Value = "'0.130"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").Value = "'138.97"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  +0.59%  "
